# Generate Report for Handback
# Replace the handback-status report data (old GUID-named files ->
# new GUID-named files) across the Overview / zh-cn / de-de sheets,
# including the filename hyperlinks' displayed text.

$wb = $excel.ActiveWorkbook

# ---- old / new identifiers -------------------------------------------------
$oldFile1 = "c701babf-ce22-4f8f-950f-b7fa872a7f71.md"
$newFile1 = "d1671c77-520e-4ab2-98a3-eed471feb362.md"

$oldFile2 = "eed7b301-1b14-4413-ba50-01a729aad8a4.md"
$newFile2 = "ffff88617cf5-5bee-48f6-a94a-84f11527f3e9.md"

$oldXlfHash = "5d44f94b859a3b2ba09089e24fa6188f32736a8a"
$newXlfHash = "82b700edb4aeb1bf4ff78bd011cb8d1ee4266f9a"

$oldXlfBase2 = "eed7b301-1b14-4413-ba50-01a729aad8a4.bca3d8c3324f5776e9a5403bfc6b688f4fbfe8a1"

$newXlfZh = "$newFile1.$newXlfHash.zh-cn.xlf" -replace "\.md", ""
$newXlfDe = "$newFile1.$newXlfHash.de-de.xlf" -replace "\.md", ""

$newXlfZh = ($newFile1 -replace "\.md$", "") + "." + $newXlfHash + ".zh-cn.xlf"
$newXlfDe = ($newFile1 -replace "\.md$", "") + "." + $newXlfHash + ".de-de.xlf"

$oldXlfZh1 = ($oldFile1 -replace "\.md$", "") + "." + $oldXlfHash + ".zh-cn.xlf"
$oldXlfDe1 = ($oldFile1 -replace "\.md$", "") + "." + $oldXlfHash + ".de-de.xlf"
$oldXlfZh2 = $oldXlfBase2 + ".zh-cn.xlf"
$oldXlfDe2 = $oldXlfBase2 + ".de-de.xlf"

$newHO = "2016-08-23 21:06:29"
$newZhGen = "2016-08-23 21:06:24"
$newZhBack = "2016-08-23 21:06:41"
$newDeGen = "2016-08-23 21:06:29"
$newDeBack = "2016-08-23 21:06:49"

# Original hyperlink target URLs (unchanged by this edit -- only the
# visible "display" text is refreshed).
$overviewUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$oldFile1"
$overviewUrl3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$oldFile2"

$zhUrlA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$oldFile1"
$zhUrlI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7dff0e96a0d8694d46ee4d7e8fc5c93045b630e1/e2e/$oldFile1"
$zhUrlA3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$oldFile2"
$zhUrlI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7dff0e96a0d8694d46ee4d7e8fc5c93045b630e1/e2e/$oldFile2"

$deUrlA2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$oldFile1"
$deUrlI2 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0e28c12bab32c45b0ac6cb0cf042504bbae55924/e2e/$oldFile1"
$deUrlA3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$oldFile2"
$deUrlI3 = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0e28c12bab32c45b0ac6cb0cf042504bbae55924/e2e/$oldFile2"

# ========================== Overview sheet ==================================
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A2").Value = $newFile1
$wsOv.Range("B2").Value = "e2e\" + $newFile1
$wsOv.Range("G2").Value = $newHO

$wsOv.Range("A3").Value = $newFile2
$wsOv.Range("B3").Value = "e2e\" + $newFile2
$wsOv.Range("G3").Value = $newHO

$wsOv.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), $overviewUrl2, "", "", "e2e\" + $newFile1)
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $overviewUrl3, "", "", "e2e\" + $newFile2)

# ============================ zh-cn sheet ===================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("G2").Value = $newXlfZh
$wsZh.Range("H2").Value = $newZhGen
$wsZh.Range("I2").Value = $newFile1
$wsZh.Range("J2").Value = $newXlfZh
$wsZh.Range("K2").Value = $newZhBack

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("G3").Value = $newXlfZh
$wsZh.Range("H3").Value = $newZhGen
$wsZh.Range("I3").Value = $newFile2
$wsZh.Range("J3").Value = $newXlfZh
$wsZh.Range("K3").Value = $newZhBack

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhUrlA2, "", "", $newFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhUrlI2, "", "", $newFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhUrlA3, "", "", $newFile2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhUrlI3, "", "", $newFile2)

# ============================ de-de sheet ===================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newDeGen
$wsDe.Range("I2").Value = $newFile1
$wsDe.Range("J2").Value = $newXlfDe
$wsDe.Range("K2").Value = $newDeBack

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newDeGen
$wsDe.Range("I3").Value = $newFile2
$wsDe.Range("J3").Value = $newXlfDe
$wsDe.Range("K3").Value = $newDeBack

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deUrlA2, "", "", $newFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deUrlI2, "", "", $newFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deUrlA3, "", "", $newFile2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deUrlI3, "", "", $newFile2)
